# Update the generated "within 100" arithmetic answer table to the
# newly generated problem set (commit: "Update master to output
# generated at 596fc94"). Each table cell holds one self-contained
# "<op>=<result>" string in its own run, so a plain whole-word
# Find/Replace keyed on the old text is unambiguous and order-independent.
$d = $word.ActiveDocument

$d.Content.Find.Execute("70-14=56", $true, $false, $false, $false, $false, $true, 1, $false, "1+17=18", 2) | Out-Null
$d.Content.Find.Execute("76-23=53", $true, $false, $false, $false, $false, $true, 1, $false, "20-17=3", 2) | Out-Null
$d.Content.Find.Execute("46+11=57", $true, $false, $false, $false, $false, $true, 1, $false, "54+36=90", 2) | Out-Null
$d.Content.Find.Execute("85+5=90", $true, $false, $false, $false, $false, $true, 1, $false, "45-33=12", 2) | Out-Null
$d.Content.Find.Execute("7+34=41", $true, $false, $false, $false, $false, $true, 1, $false, "89-23=66", 2) | Out-Null
$d.Content.Find.Execute("84-83=1", $true, $false, $false, $false, $false, $true, 1, $false, "48-19=29", 2) | Out-Null
$d.Content.Find.Execute("29+69=98", $true, $false, $false, $false, $false, $true, 1, $false, "80-1=79", 2) | Out-Null
$d.Content.Find.Execute("8+70=78", $true, $false, $false, $false, $false, $true, 1, $false, "31-16=15", 2) | Out-Null
$d.Content.Find.Execute("44+15=59", $true, $false, $false, $false, $false, $true, 1, $false, "5+8=13", 2) | Out-Null
$d.Content.Find.Execute("86-80=6", $true, $false, $false, $false, $false, $true, 1, $false, "11+83=94", 2) | Out-Null
$d.Content.Find.Execute("76-42=34", $true, $false, $false, $false, $false, $true, 1, $false, "25+48=73", 2) | Out-Null
$d.Content.Find.Execute("46-11=35", $true, $false, $false, $false, $false, $true, 1, $false, "25+2=27", 2) | Out-Null
$d.Content.Find.Execute("55+26=81", $true, $false, $false, $false, $false, $true, 1, $false, "55+35=90", 2) | Out-Null
$d.Content.Find.Execute("58-56=2", $true, $false, $false, $false, $false, $true, 1, $false, "60-55=5", 2) | Out-Null
$d.Content.Find.Execute("16+57=73", $true, $false, $false, $false, $false, $true, 1, $false, "59-6=53", 2) | Out-Null
$d.Content.Find.Execute("65+4=69", $true, $false, $false, $false, $false, $true, 1, $false, "20+11=31", 2) | Out-Null
$d.Content.Find.Execute("61-19=42", $true, $false, $false, $false, $false, $true, 1, $false, "59-10=49", 2) | Out-Null
$d.Content.Find.Execute("67-41=26", $true, $false, $false, $false, $false, $true, 1, $false, "62+36=98", 2) | Out-Null
$d.Content.Find.Execute("55-45=10", $true, $false, $false, $false, $false, $true, 1, $false, "62-47=15", 2) | Out-Null
$d.Content.Find.Execute("81-74=7", $true, $false, $false, $false, $false, $true, 1, $false, "81-61=20", 2) | Out-Null
$d.Content.Find.Execute("35-32=3", $true, $false, $false, $false, $false, $true, 1, $false, "87-28=59", 2) | Out-Null
$d.Content.Find.Execute("56-13=43", $true, $false, $false, $false, $false, $true, 1, $false, "69-58=11", 2) | Out-Null
$d.Content.Find.Execute("31-3=28", $true, $false, $false, $false, $false, $true, 1, $false, "19+30=49", 2) | Out-Null
$d.Content.Find.Execute("49-20=29", $true, $false, $false, $false, $false, $true, 1, $false, "36+34=70", 2) | Out-Null
$d.Content.Find.Execute("30+66=96", $true, $false, $false, $false, $false, $true, 1, $false, "39-2=37", 2) | Out-Null
$d.Content.Find.Execute("95-67=28", $true, $false, $false, $false, $false, $true, 1, $false, "68-30=38", 2) | Out-Null
$d.Content.Find.Execute("96-20=76", $true, $false, $false, $false, $false, $true, 1, $false, "51+8=59", 2) | Out-Null
$d.Content.Find.Execute("42+19=61", $true, $false, $false, $false, $false, $true, 1, $false, "84-53=31", 2) | Out-Null
$d.Content.Find.Execute("58-26=32", $true, $false, $false, $false, $false, $true, 1, $false, "4+43=47", 2) | Out-Null
$d.Content.Find.Execute("5+33=38", $true, $false, $false, $false, $false, $true, 1, $false, "68+16=84", 2) | Out-Null
$d.Content.Find.Execute("72-32=40", $true, $false, $false, $false, $false, $true, 1, $false, "74-0=74", 2) | Out-Null
$d.Content.Find.Execute("83-73=10", $true, $false, $false, $false, $false, $true, 1, $false, "55+11=66", 2) | Out-Null
$d.Content.Find.Execute("99-88=11", $true, $false, $false, $false, $false, $true, 1, $false, "63-37=26", 2) | Out-Null
$d.Content.Find.Execute("31+34=65", $true, $false, $false, $false, $false, $true, 1, $false, "98-55=43", 2) | Out-Null
$d.Content.Find.Execute("22+57=79", $true, $false, $false, $false, $false, $true, 1, $false, "11+2=13", 2) | Out-Null
$d.Content.Find.Execute("38-30=8", $true, $false, $false, $false, $false, $true, 1, $false, "42-31=11", 2) | Out-Null
$d.Content.Find.Execute("9+73=82", $true, $false, $false, $false, $false, $true, 1, $false, "78+2=80", 2) | Out-Null
$d.Content.Find.Execute("60+34=94", $true, $false, $false, $false, $false, $true, 1, $false, "43-18=25", 2) | Out-Null
$d.Content.Find.Execute("7+43=50", $true, $false, $false, $false, $false, $true, 1, $false, "15+78=93", 2) | Out-Null
$d.Content.Find.Execute("43-4=39", $true, $false, $false, $false, $false, $true, 1, $false, "55+24=79", 2) | Out-Null
$d.Content.Find.Execute("23+68=91", $true, $false, $false, $false, $false, $true, 1, $false, "64-38=26", 2) | Out-Null
$d.Content.Find.Execute("82-27=55", $true, $false, $false, $false, $false, $true, 1, $false, "76-3=73", 2) | Out-Null
$d.Content.Find.Execute("85-42=43", $true, $false, $false, $false, $false, $true, 1, $false, "3+69=72", 2) | Out-Null
$d.Content.Find.Execute("73-35=38", $true, $false, $false, $false, $false, $true, 1, $false, "8-1=7", 2) | Out-Null
$d.Content.Find.Execute("14+57=71", $true, $false, $false, $false, $false, $true, 1, $false, "58+24=82", 2) | Out-Null
$d.Content.Find.Execute("84-74=10", $true, $false, $false, $false, $false, $true, 1, $false, "36+38=74", 2) | Out-Null
$d.Content.Find.Execute("11+14=25", $true, $false, $false, $false, $false, $true, 1, $false, "17+32=49", 2) | Out-Null
$d.Content.Find.Execute("84+3=87", $true, $false, $false, $false, $false, $true, 1, $false, "85-47=38", 2) | Out-Null
$d.Content.Find.Execute("8+55=63", $true, $false, $false, $false, $false, $true, 1, $false, "1+12=13", 2) | Out-Null
$d.Content.Find.Execute("11+59=70", $true, $false, $false, $false, $false, $true, 1, $false, "38+44=82", 2) | Out-Null
$d.Content.Find.Execute("49+3=52", $true, $false, $false, $false, $false, $true, 1, $false, "87+7=94", 2) | Out-Null
$d.Content.Find.Execute("24-1=23", $true, $false, $false, $false, $false, $true, 1, $false, "55+24=79", 2) | Out-Null
$d.Content.Find.Execute("83-61=22", $true, $false, $false, $false, $false, $true, 1, $false, "72-39=33", 2) | Out-Null
$d.Content.Find.Execute("86-46=40", $true, $false, $false, $false, $false, $true, 1, $false, "80-2=78", 2) | Out-Null
$d.Content.Find.Execute("30-19=11", $true, $false, $false, $false, $false, $true, 1, $false, "88-48=40", 2) | Out-Null
$d.Content.Find.Execute("0+15=15", $true, $false, $false, $false, $false, $true, 1, $false, "13+30=43", 2) | Out-Null
$d.Content.Find.Execute("33+52=85", $true, $false, $false, $false, $false, $true, 1, $false, "17+30=47", 2) | Out-Null
$d.Content.Find.Execute("17+56=73", $true, $false, $false, $false, $false, $true, 1, $false, "35-21=14", 2) | Out-Null
$d.Content.Find.Execute("35+12=47", $true, $false, $false, $false, $false, $true, 1, $false, "15-9=6", 2) | Out-Null
$d.Content.Find.Execute("21+67=88", $true, $false, $false, $false, $false, $true, 1, $false, "73-54=19", 2) | Out-Null
$d.Content.Find.Execute("48+26=74", $true, $false, $false, $false, $false, $true, 1, $false, "73-51=22", 2) | Out-Null
$d.Content.Find.Execute("61+4=65", $true, $false, $false, $false, $false, $true, 1, $false, "12+13=25", 2) | Out-Null
$d.Content.Find.Execute("90-66=24", $true, $false, $false, $false, $false, $true, 1, $false, "61-41=20", 2) | Out-Null
$d.Content.Find.Execute("35+4=39", $true, $false, $false, $false, $false, $true, 1, $false, "2+42=44", 2) | Out-Null
$d.Content.Find.Execute("24+39=63", $true, $false, $false, $false, $false, $true, 1, $false, "83-10=73", 2) | Out-Null
$d.Content.Find.Execute("13+59=72", $true, $false, $false, $false, $false, $true, 1, $false, "33+34=67", 2) | Out-Null
$d.Content.Find.Execute("9+42=51", $true, $false, $false, $false, $false, $true, 1, $false, "9+66=75", 2) | Out-Null
$d.Content.Find.Execute("44+16=60", $true, $false, $false, $false, $false, $true, 1, $false, "26+70=96", 2) | Out-Null
$d.Content.Find.Execute("44+39=83", $true, $false, $false, $false, $false, $true, 1, $false, "24+18=42", 2) | Out-Null
$d.Content.Find.Execute("87+9=96", $true, $false, $false, $false, $false, $true, 1, $false, "74-69=5", 2) | Out-Null
$d.Content.Find.Execute("88+7=95", $true, $false, $false, $false, $false, $true, 1, $false, "74-13=61", 2) | Out-Null
$d.Content.Find.Execute("12+40=52", $true, $false, $false, $false, $false, $true, 1, $false, "0+7=7", 2) | Out-Null
$d.Content.Find.Execute("59-5=54", $true, $false, $false, $false, $false, $true, 1, $false, "43-22=21", 2) | Out-Null
$d.Content.Find.Execute("34+45=79", $true, $false, $false, $false, $false, $true, 1, $false, "83+13=96", 2) | Out-Null
$d.Content.Find.Execute("9+10=19", $true, $false, $false, $false, $false, $true, 1, $false, "11+39=50", 2) | Out-Null
$d.Content.Find.Execute("90-62=28", $true, $false, $false, $false, $false, $true, 1, $false, "39-37=2", 2) | Out-Null
$d.Content.Find.Execute("2+31=33", $true, $false, $false, $false, $false, $true, 1, $false, "15+52=67", 2) | Out-Null
$d.Content.Find.Execute("73-8=65", $true, $false, $false, $false, $false, $true, 1, $false, "96-60=36", 2) | Out-Null
$d.Content.Find.Execute("73+25=98", $true, $false, $false, $false, $false, $true, 1, $false, "10-1=9", 2) | Out-Null
$d.Content.Find.Execute("50-31=19", $true, $false, $false, $false, $false, $true, 1, $false, "49-33=16", 2) | Out-Null
$d.Content.Find.Execute("69-17=52", $true, $false, $false, $false, $false, $true, 1, $false, "60-31=29", 2) | Out-Null
$d.Content.Find.Execute("83-37=46", $true, $false, $false, $false, $false, $true, 1, $false, "51+17=68", 2) | Out-Null
$d.Content.Find.Execute("29+26=55", $true, $false, $false, $false, $false, $true, 1, $false, "89+1=90", 2) | Out-Null
$d.Content.Find.Execute("27+27=54", $true, $false, $false, $false, $false, $true, 1, $false, "9-1=8", 2) | Out-Null
$d.Content.Find.Execute("13-2=11", $true, $false, $false, $false, $false, $true, 1, $false, "95-9=86", 2) | Out-Null
$d.Content.Find.Execute("89-59=30", $true, $false, $false, $false, $false, $true, 1, $false, "34-8=26", 2) | Out-Null
$d.Content.Find.Execute("49-37=12", $true, $false, $false, $false, $false, $true, 1, $false, "78+13=91", 2) | Out-Null
$d.Content.Find.Execute("47+37=84", $true, $false, $false, $false, $false, $true, 1, $false, "85-50=35", 2) | Out-Null
$d.Content.Find.Execute("85+3=88", $true, $false, $false, $false, $false, $true, 1, $false, "12-9=3", 2) | Out-Null
$d.Content.Find.Execute("51-6=45", $true, $false, $false, $false, $false, $true, 1, $false, "97-73=24", 2) | Out-Null
$d.Content.Find.Execute("16+23=39", $true, $false, $false, $false, $false, $true, 1, $false, "45+47=92", 2) | Out-Null
$d.Content.Find.Execute("14+40=54", $true, $false, $false, $false, $false, $true, 1, $false, "96+3=99", 2) | Out-Null
$d.Content.Find.Execute("22+58=80", $true, $false, $false, $false, $false, $true, 1, $false, "12-5=7", 2) | Out-Null
$d.Content.Find.Execute("22+34=56", $true, $false, $false, $false, $false, $true, 1, $false, "29+23=52", 2) | Out-Null
$d.Content.Find.Execute("46+3=49", $true, $false, $false, $false, $false, $true, 1, $false, "70+27=97", 2) | Out-Null
$d.Content.Find.Execute("53+9=62", $true, $false, $false, $false, $false, $true, 1, $false, "36+32=68", 2) | Out-Null
$d.Content.Find.Execute("24+27=51", $true, $false, $false, $false, $false, $true, 1, $false, "72-5=67", 2) | Out-Null
$d.Content.Find.Execute("67-16=51", $true, $false, $false, $false, $false, $true, 1, $false, "96-36=60", 2) | Out-Null
$d.Content.Find.Execute("41+26=67", $true, $false, $false, $false, $false, $true, 1, $false, "58+4=62", 2) | Out-Null
$d.Content.Find.Execute("60-49=11", $true, $false, $false, $false, $false, $true, 1, $false, "90-45=45", 2) | Out-Null
